$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking.com crypto price/volume refresh (GitHub Actions bot).
# Price/Volume columns are stored as text in this sheet. Force plain-text
# format on any Price cell whose new value would otherwise be read by Excel
# as a number (e.g. '1.001', '17.00') so trailing zeros and the original
# text representation are preserved, matching the source data feed.

$textFormatCells = @('D4', 'D5', 'D6', 'D7', 'D9', 'D11', 'D12', 'D13', 'D14', 'D15', 'D18', 'D19', 'D20', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D30', 'D31', 'D32', 'D33', 'D34', 'D37', 'D38', 'D39', 'D40', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '28.139.18'
$ws.Range('E2').Value = '  -0.14%  '

# Row 3
$ws.Range('D3').Value = '1.761.01'
$ws.Range('E3').Value = '  -2.76%  '

# Row 4
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.29%  '

# Row 5
$ws.Range('D5').Value = '334.96'
$ws.Range('E5').Value = '  -1.12%  '

# Row 6
$ws.Range('D6').Value = '0.9978'
$ws.Range('E6').Value = '  -0.17%  '

# Row 7
$ws.Range('D7').Value = '0.3780'
$ws.Range('E7').Value = '  -3.23%  '

# Row 8
$ws.Range('E8').Value = '  -3.32%  '

# Row 9
$ws.Range('D9').Value = '45.61'
$ws.Range('E9').Value = '  -5.80%  '

# Row 10
$ws.Range('E10').Value = '  -5.26%  '

# Row 11
$ws.Range('D11').Value = '0.07202'
$ws.Range('E11').Value = '  -4.57%  '

# Row 12
$ws.Range('B12').Value = 'BinanceUSD'
$ws.Range('C12').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D12').Value = '0.9997'
$ws.Range('E12').Value = '  -0.20%  '

# Row 13
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = '22.58'
$ws.Range('E13').Value = '  +2.49%  '

# Row 14
$ws.Range('D14').Value = '6.187'
$ws.Range('E14').Value = '  -4.88%  '

# Row 15
$ws.Range('D15').Value = '7.217'
$ws.Range('E15').Value = '  +1.19%  '

# Row 16
$ws.Range('D16').Value = '1.758.09'
$ws.Range('E16').Value = '  -3.14%  '

# Row 17
$ws.Range('E17').Value = '  -4.13%  '

# Row 18
$ws.Range('D18').Value = '0.06575'
$ws.Range('E18').Value = '  -1.69%  '

# Row 19
$ws.Range('D19').Value = '80.89'
$ws.Range('E19').Value = '  -4.73%  '

# Row 20
$ws.Range('D20').Value = '0.9979'

# Row 21
$ws.Range('D21').Value = '17.00'
$ws.Range('E21').Value = '  -4.01%  '

# Row 22
$ws.Range('D22').Value = '6.273'
$ws.Range('E22').Value = '  -4.19%  '

# Row 23
$ws.Range('D23').Value = '28.120.50'
$ws.Range('E23').Value = '  -0.25%  '

# Row 24
$ws.Range('D24').Value = '11.64'
$ws.Range('E24').Value = '  -6.35%  '

# Row 25
$ws.Range('D25').Value = '2.385'
$ws.Range('E25').Value = '  -0.84%  '

# Row 26
$ws.Range('D26').Value = '153.39'
$ws.Range('E26').Value = '  -0.28%  '

# Row 27
$ws.Range('D27').Value = '19.87'
$ws.Range('E27').Value = '  -6.58%  '

# Row 28
$ws.Range('D28').Value = '2.333'
$ws.Range('E28').Value = '  -7.53%  '

# Row 29
$ws.Range('D29').Value = '1.958.95'
$ws.Range('E29').Value = '  -3.03%  '

# Row 30
$ws.Range('D30').Value = '1.270'
$ws.Range('E30').Value = '  -14.46%  '

# Row 31
$ws.Range('D31').Value = '131.93'
$ws.Range('E31').Value = '  -2.68%  '

# Row 32
$ws.Range('D32').Value = '4.015'
$ws.Range('E32').Value = '  +0.03%  '

# Row 33
$ws.Range('D33').Value = '5.814'
$ws.Range('E33').Value = '  -5.28%  '

# Row 34
$ws.Range('D34').Value = '0.08793'
$ws.Range('E34').Value = '  -0.29%  '

# Row 35
$ws.Range('E35').Value = '  -5.62%  '

# Row 36
$ws.Range('E36').Value = '  -3.04%  '

# Row 37
$ws.Range('D37').Value = '0.6632'
$ws.Range('E37').Value = '  -4.43%  '

# Row 38
$ws.Range('D38').Value = '0.06206'
$ws.Range('E38').Value = '  -5.10%  '

# Row 39
$ws.Range('D39').Value = '5.174'
$ws.Range('E39').Value = '  -5.11%  '

# Row 40
$ws.Range('D40').Value = '0.2116'
$ws.Range('E40').Value = '  -4.10%  '

# Row 41
$ws.Range('E41').Value = '  -2.86%  '

# Row 42
$ws.Range('D42').Value = '1.449'
$ws.Range('E42').Value = '  -9.90%  '

# Row 43
$ws.Range('D43').Value = '8.027'
$ws.Range('E43').Value = '  -4.88%  '

# Row 44
$ws.Range('D44').Value = '0.9976'
$ws.Range('E44').Value = '  -0.11%  '

# Row 45
$ws.Range('D45').Value = '13.77'
$ws.Range('E45').Value = '  -5.69%  '

# Row 46
$ws.Range('D46').Value = '0.6048'
$ws.Range('E46').Value = '  -5.89%  '

# Row 47
$ws.Range('D47').Value = '3.815'
$ws.Range('E47').Value = '  -1.25%  '

# Row 48
$ws.Range('D48').Value = '129.56'
$ws.Range('E48').Value = '  -1.47%  '

# Row 49
$ws.Range('D49').Value = '2.014'
$ws.Range('E49').Value = '  -6.03%  '

# Row 50
$ws.Range('D50').Value = '1.189'
$ws.Range('E50').Value = '  +3.18%  '

# Row 51
$ws.Range('D51').Value = '0.07210'
$ws.Range('E51').Value = '  +0.19%  '
